$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 corresponds to the "CLOPEX GREL 75MG 30 F.C.TAB" line item which has
# been removed from this day's sale report. Deleting the entire row shifts
# every row below it up by one (renumbering handled automatically for the
# running index column) and Excel will compact the shared-strings table and
# merged-cell list on save.
$ws.Rows.Item(20).Delete()

# Update the "generated at" timestamp footer (now on row 121 after the row
# delete above) to reflect the new export time.
$ws.Range("A121").Value = "Monday, 21 July, 2025 8:18 PM"
